# Saldo.xlsx update:
#  - The account 004001621 / DANIELA row is moved from its old position
#    (row 7, right after LUCIA) down to just after the PEDRO/834.33 row,
#    and its balance (column C) changes from 29833.09 to 833.09.
#  - The account 004448303 / NASSIM row (balance 7569.82) is removed
#    entirely.
#  - All other rows keep their relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the old DANIELA (004001621) row at row 7. Everything below
#    shifts up by one row.
$ws.Rows(7).Delete()

# 2) The NASSIM (004448303) row, originally row 11, is now at row 10
#    after the shift above. Remove it entirely.
$ws.Rows(10).Delete()

# 3) Insert a new row right after the PEDRO/834.33 row (now row 24) to
#    hold the relocated DANIELA entry with its updated balance.
$ws.Rows(25).Insert()

# Force column A to be treated as text so the leading zeros in the
# account number are preserved, matching the rest of the sheet.
$ws.Cells.Item(25, 1).NumberFormat = "@"
$ws.Cells.Item(25, 1).Value = "004001621"
$ws.Cells.Item(25, 2).Value = "DANIELA"
$ws.Cells.Item(25, 3).Value = 833.09
